# Add a new data row (55) to the PEBCOM sheet, following the same
# pattern as every existing row: columns A-L are stored as text
# (even when the content looks numeric, e.g. case numbers, OT numbers,
# comuna codes, attachment counts) and columns M-N (coordinates) are
# stored as real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 55

# Force columns A:L to Text format *before* writing, otherwise values
# like "-483", "13", "807605744", "1" or the "6/18/2025" date string
# would be auto-coerced into numbers/dates by Excel's type inference.
$textRange = $ws.Range("A" + $row + ":L" + $row)
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value  = "-483"
$ws.Cells.Item($row, 2).Value  = "6/18/2025"
$ws.Cells.Item($row, 3).Value  = "Arcos 4326"
$ws.Cells.Item($row, 4).Value  = "13"
$ws.Cells.Item($row, 5).Value  = "807605744"
$ws.Cells.Item($row, 6).Value  = "PEBCOM"
$ws.Cells.Item($row, 7).Value  = "Pendiente"
$ws.Cells.Item($row, 8).Value  = "Terminal de teco con clientes Se solicita desconectar alarma vecinal en el poste"
$ws.Cells.Item($row, 9).Value  = "1"
$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Sin equipos"
$ws.Cells.Item($row, 12).Value = "Poste"

# Coordinates stay numeric.
$ws.Cells.Item($row, 13).Value = -58.469257
$ws.Cells.Item($row, 14).Value = -34.542018
